# Updated cryptos list (prices/volume) per the Thu Sep 28 11:24:14 UTC 2023
# GitHub Actions refresh. Price/volume cells hold plain text (e.g. "1.00",
# "0.999") so numeric-looking values are written via FormulaR1C1 with a
# leading apostrophe to keep them as text instead of being parsed as
# numbers (which would silently drop formatting like the trailing zero in
# "1.00"). Non-numeric-looking text (names, URLs, percent strings) is set
# directly through .Value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.494.05"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.625.60"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").FormulaR1C1 = "'0.999"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").FormulaR1C1 = "'213.51"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").FormulaR1C1 = "'0.502"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").FormulaR1C1 = "'1.00"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").FormulaR1C1 = "'19.19"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").FormulaR1C1 = "'0.0853"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.854.58"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "1.630.23"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").FormulaR1C1 = "'0.511"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").FormulaR1C1 = "'64.02"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").FormulaR1C1 = "'234.99"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "26.489.83"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").FormulaR1C1 = "'7.78"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").FormulaR1C1 = "'1.00"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").FormulaR1C1 = "'4.31"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").FormulaR1C1 = "'2.20"
$ws.Range("E23").Value = "  +2.98%  "
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").FormulaR1C1 = "'9.14"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").FormulaR1C1 = "'147.19"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").FormulaR1C1 = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").FormulaR1C1 = "'7.07"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").FormulaR1C1 = "'15.66"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "1.522.82"
$ws.Range("E32").Value = "  +4.75%  "
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").FormulaR1C1 = "'2.98"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").FormulaR1C1 = "'0.568"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").FormulaR1C1 = "'0.0166"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").FormulaR1C1 = "'0.835"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").FormulaR1C1 = "'5.86"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").FormulaR1C1 = "'1.00"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "1.765.30"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").FormulaR1C1 = "'63.05"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").FormulaR1C1 = "'90.13"
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("D48").FormulaR1C1 = "'1.51"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").FormulaR1C1 = "'0.0964"
$ws.Range("E51").Value = "  +0.15%  "
